$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 444, shifting existing rows 444:507 down to 445:508.
$ws.Rows(444).Insert()

# Populate the newly inserted row 444 with its data (matches the surrounding
# rows for the constant columns, and new values for the varying columns).
$ws.Cells.Item(444, 1).Value = 3
$ws.Cells.Item(444, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(444, 3).Value = "Coquimbo"
$ws.Cells.Item(444, 4).Value = 45154
$ws.Cells.Item(444, 5).Value = 5
$ws.Cells.Item(444, 6).Value = 100112001
$ws.Cells.Item(444, 7).Value = "Berenjena"
$ws.Cells.Item(444, 8).Value = "Sin especificar"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 50
$ws.Cells.Item(444, 11).Value = 7500
$ws.Cells.Item(444, 12).Value = 7500
$ws.Cells.Item(444, 13).Value = 7500
$ws.Cells.Item(444, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(444, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(444, 16).Value = 125
$ws.Cells.Item(444, 17).Value = 60
$ws.Cells.Item(444, 18).Value = "Hortaliza"
